# Adds one new weekly price observation (2 quality-grade rows: "Primera"
# and "Segunda") for "Crespo record" Repollo at Terminal Hortofrutícola
# Agro Chillán, inserted right after the existing row 439 (row 440 in the
# resulting sheet), which pushes every following row down by two and
# grows the used range from A1:R524 to A1:R526.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 440/441; everything that used to live at row
# 440 onward shifts down to 442 onward.
$ws.Range("A440:R441").Insert()

# New row 440: "Primera" quality, $1.000 flat, volume 300.
$ws.Range("A440").Value = 7
$ws.Range("B440").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C440").Value = "Ñuble"
$ws.Range("D440").Value = 45209
$ws.Range("E440").Value = 16
$ws.Range("F440").Value = 100112006
$ws.Range("G440").Value = "Repollo"
$ws.Range("H440").Value = "Crespo record"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 300
$ws.Range("K440").Value = 1000
$ws.Range("L440").Value = 1000
$ws.Range("M440").Value = 1000
$ws.Range("N440").Value = "$/unidad"
$ws.Range("O440").Value = "Provincia de Diguillín"
$ws.Range("P440").Value = 1000
$ws.Range("Q440").Value = 1
$ws.Range("R440").Value = "Hortaliza"

# New row 441: "Segunda" quality, $800 flat, volume 300.
$ws.Range("A441").Value = 7
$ws.Range("B441").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C441").Value = "Ñuble"
$ws.Range("D441").Value = 45209
$ws.Range("E441").Value = 16
$ws.Range("F441").Value = 100112006
$ws.Range("G441").Value = "Repollo"
$ws.Range("H441").Value = "Crespo record"
$ws.Range("I441").Value = "Segunda"
$ws.Range("J441").Value = 300
$ws.Range("K441").Value = 800
$ws.Range("L441").Value = 800
$ws.Range("M441").Value = 800
$ws.Range("N441").Value = "$/unidad"
$ws.Range("O441").Value = "Provincia de Diguillín"
$ws.Range("P441").Value = 800
$ws.Range("Q441").Value = 1
$ws.Range("R441").Value = "Hortaliza"
